$d = $word.ActiveDocument

$d.Content.Find.Execute("traoudant no where", $true, $false, $false, $false, $false,
                         $true, 1, $false, "traoudant we are from", 2)

$d.Content.Find.Execute("       RR000000019MA", $true, $false, $false, $false, $false,
                         $true, 1, $false, "       RR000000022MA", 2)

$d.Content.Find.Execute("MR.qwqw qwqw", $true, $false, $false, $false, $false,
                         $true, 1, $false, "MR.qw qw", 2)

$d.Content.Find.Execute("Ain Chegga : erqw", $true, $false, $false, $false, $false,
                         $true, 1, $false, "Ain Chegga : qw", 2)

$d.Content.Find.Execute("123123", $true, $false, $false, $false, $false,
                         $true, 1, $false, "12331", 2)
